# Add a new "2020" column (column N) mirroring column M's styling, with
# new data values, then adjust the view/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + data values for column N (2020), one per row 4-17.
$ws.Range("N4").Value = 2020
$ws.Range("N5").Value = 11.4
$ws.Range("N6").Value = 14.7
$ws.Range("N7").Value = 9
$ws.Range("N8").Value = 10.8
$ws.Range("N9").Value = 4.7
$ws.Range("N10").Value = 5.0999999999999996
$ws.Range("N11").Value = 3.4
$ws.Range("N12").Value = 19.7
$ws.Range("N13").Value = 18.8
$ws.Range("N14").Value = 6.8
$ws.Range("N16").Value = 12.5
$ws.Range("N17").Value = 10.7

# Copy column M's formatting onto column N for the rows that have data so
# the new column visually matches the rest of the table (N15 stays blank,
# matching M15).
$ws.Range("M4:M17").Copy()
$ws.Range("N4:N17").PasteSpecial(-4122) # xlPasteFormats

# Update the view to scroll to the new column and move the selection.
$ws.Application.ActiveWindow.ScrollColumn = 5
$sheetView = $ws.Application.ActiveWindow
$sheetView.TopLeftCell = $ws.Range("E2")
$ws.Range("S18").Select()
